$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes -----------------------------------------------
# Drop the old data row 4 (Site Ref 61712824) entirely - it was removed in
# the "redone boundary" rework, shifting all the rows below it up by one.
$ws.Rows.Item(4).EntireRow.Delete() | Out-Null

# The two bottom-most rows (previously rows 15 and 16, now rows 14 and 15
# after the shift above) were also dropped, leaving 13 rows total.
$ws.Rows.Item(14).EntireRow.Delete() | Out-Null
$ws.Rows.Item(14).EntireRow.Delete() | Out-Null

# The old "Season" column (O) was removed; the "Sample timeframe" column
# (old P) slides left to become the new column O.
$ws.Columns.Item(15).EntireColumn.Delete() | Out-Null

# --- Corrected values for the re-derived boundary rows -----------------
# (new rows 8-12, i.e. the old AR7/AR9/AR12(ART)/AR15/AR16 records)
$ws.Range("B8").Value = 25408
$ws.Range("C8").Value = 7
$ws.Range("F8").Value = 8.51
$ws.Range("G8").Value = 214.604
$ws.Range("N8").Value = 215.8364202880859

$ws.Range("B9").Value = 25408
$ws.Range("C9").Value = 7
$ws.Range("F9").Value = 15.83
$ws.Range("G9").Value = 215.419
$ws.Range("N9").Value = 217.9015521240234

$ws.Range("B10").Value = 25373
$ws.Range("C10").Value = 6
$ws.Range("F10").Value = -2.52
$ws.Range("G10").Value = 211.802
$ws.Range("N10").Value = 213.3845935058594

$ws.Range("B11").Value = 25408
$ws.Range("C11").Value = 7
$ws.Range("F11").Value = 9.02
$ws.Range("G11").Value = 213.835
$ws.Range("N11").Value = 214.0011029052734

$ws.Range("B12").Value = 25408
$ws.Range("C12").Value = 7
$ws.Range("F12").Value = 59.29
$ws.Range("G12").Value = 212.97
$ws.Range("N12").Value = 212.7830590820312

Write-Output "edit applied"
